$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the existing date-formatted style (used by A2/A3) to the new date cells,
# so the new cells reuse the same style index instead of creating new ones.
$ws.Cells.Item(2, 1).Copy()
$ws.Range("A4:A16").PasteSpecial(-4122)

# Add new rows of data (dates in column A, zeros in column B) for rows 4-16
$dates = 45939,45940,45941,45942,45943,45944,45945,45946,45947,45948,45949,45950,45951
$row = 4
foreach ($d in $dates) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 2).Value = 0
    $row = $row + 1
}

# Update the active selection to match the target state
$ws.Range("L19").Select()

# Touch the page setup so the orientation is written out explicitly
$ws.PageSetup.Orientation = 1
